# "ssh key instructions added"
#
# Adds a new "Connecting to GitHub with SSH" section to the Settings sheet,
# with the ssh-keygen / ssh-agent / clip commands and their descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- New section header (row 27) ------------------------------------------
$ws.Range("A27").Value2 = "Connecting to GitHub with SSH"
$ws.Range("A27").Font.Bold = $true

# --- Command / description rows (29-31) ------------------------------------
$ws.Range("A29").Value2 = 'ssh-keygen -o -a 100 -t ed25519 -f ~/.ssh/id_ed25519 -C "mkrstv@gmail.com"'
$ws.Range("B29").Value2 = "generates the key"

$ws.Range("A30").Value2 = "eval ``ssh-agent -s``"
$ws.Range("B30").Value2 = "start the ssh-agent in the background"

# Extra note columns added alongside the ssh-keygen row
$ws.Range("C29").Value2 = "passphrase"
$ws.Range("D29").Value2 = "LeeMu"

# Row 31 previously held an empty spacer cell (A31) with a different font
# size/style; reset its formatting to the normal body style before writing
# into it so it matches the other data rows (copy the already-normalized
# format from A29/B29).
$ws.Range("A29:B29").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A31").Value2 = "clip < ~/.ssh/id_ed25519.pub"
$ws.Range("B31").Value2 = "Copy the SSH public key to your clipboard"

# Leave the selection on the new section header, matching where the user
# ended up after typing in the new content.
$ws.Activate() | Out-Null
$ws.Range("A27").Select() | Out-Null
